$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 7013
$ws.Range("I20").Value = 1684
$ws.Range("J20").Value = 23000
$ws.Range("K20").Value = 1684
$ws.Range("L20").Value = 23000
$ws.Range("M20").Value = -1454
$ws.Range("N20").Value = -23460

$ws.Range("H35").Value = 7013
$ws.Range("I35").Value = 1684
$ws.Range("J35").Value = 23000
$ws.Range("K35").Value = 1684
$ws.Range("L35").Value = 23000
$ws.Range("M35").Value = -1305
$ws.Range("N35").Value = -23758

$ws.Range("H116").Value = 3285.5334
$ws.Range("I116").Value = 2790.4443
$ws.Range("J116").Value = 4028.1667
$ws.Range("K116").Value = 2790.4443
$ws.Range("L116").Value = 4028.1667
$ws.Range("M116").Value = 651.5556999999999
$ws.Range("N116").Value = -10912.1667

$ws.Range("H125").Value = 667.0909
$ws.Range("I125").Value = 637.8125
$ws.Range("J125").Value = 745.1667
$ws.Range("K125").Value = 5740.3125
$ws.Range("L125").Value = 6706.5003
$ws.Range("M125").Value = -3280.3125
$ws.Range("N125").Value = -11626.5003

$ws.Range("H132").Value = 5866.25
$ws.Range("I132").Value = 2132.838
$ws.Range("J132").Value = 25600
$ws.Range("K132").Value = 6398.514000000001
$ws.Range("L132").Value = 76800
$ws.Range("M132").Value = -3868.514000000001
$ws.Range("N132").Value = -81860

$ws.Range("H135").Value = 477.33334
$ws.Range("I135").Value = 224.5
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 2020.5
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = 514.5
$ws.Range("N135").Value = -27570

$ws.Range("H136").Value = 28994.5
$ws.Range("J136").Value = 28994.5
$ws.Range("L136").Value = 28994.5
$ws.Range("N136").Value = -39194.5

$ws.Range("H140").Value = 42858
$ws.Range("J140").Value = 42858
$ws.Range("L140").Value = 42858
$ws.Range("N140").Value = -53218

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 764.2857
$ws.Range("I45").Value = 691.6667
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 691.6667
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -314.6667
$ws.Range("N45").Value = -1954

$ws.Range("H61").Value = 591066.5
$ws.Range("I61").Value = 528501.8
$ws.Range("J61").Value = 670315.2
$ws.Range("K61").Value = 528501.8
$ws.Range("L61").Value = 670315.2
$ws.Range("M61").Value = -528289.8
$ws.Range("N61").Value = -670739.2

$ws.Range("H74").Value = 11685147
$ws.Range("I74").Value = 8374612
$ws.Range("J74").Value = 19669378
$ws.Range("K74").Value = 8374612
$ws.Range("L74").Value = 19669378
$ws.Range("M74").Value = -8373738
$ws.Range("N74").Value = -19671126

$ws.Range("H77").Value = 11685147
$ws.Range("I77").Value = 8374612
$ws.Range("J77").Value = 19669378
$ws.Range("K77").Value = 41873060
$ws.Range("L77").Value = 98346890
$ws.Range("M77").Value = -41868692
$ws.Range("N77").Value = -98355626

$ws.Range("H132").Value = 2109.4468
$ws.Range("I132").Value = 1182.5385
$ws.Range("J132").Value = 3257.0476
$ws.Range("K132").Value = 3547.6155
$ws.Range("L132").Value = 9771.1428
$ws.Range("M132").Value = -1017.6155
$ws.Range("N132").Value = -14831.1428

$ws.Range("H136").Value = 591066.5
$ws.Range("I136").Value = 528501.8
$ws.Range("J136").Value = 670315.2
$ws.Range("K136").Value = 1585505.4
$ws.Range("L136").Value = 2010945.6
$ws.Range("M136").Value = -1582955.4
$ws.Range("N136").Value = -2016045.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H36").Value = 4450
$ws.Range("J36").Value = 6000
$ws.Range("L36").Value = 6000
$ws.Range("N36").Value = -7068

$ws.Range("H39").Value = 4500
$ws.Range("J39").Value = 4500
$ws.Range("L39").Value = 4500
$ws.Range("N39").Value = -5278

$ws.Range("H54").Value = 3059.4666
$ws.Range("I54").Value = 1265.091
$ws.Range("J54").Value = 7994
$ws.Range("K54").Value = 1265.091
$ws.Range("L54").Value = 7994
$ws.Range("M54").Value = -781.0909999999999
$ws.Range("N54").Value = -8962

$ws.Range("H105").Value = 1724.5454
$ws.Range("I105").Value = 1622.4642
$ws.Range("J105").Value = 2296.2
$ws.Range("K105").Value = 1622.4642
$ws.Range("L105").Value = 2296.2
$ws.Range("M105").Value = 124.5358000000001
$ws.Range("N105").Value = -5790.2

$ws.Range("H126").Value = 59000
$ws.Range("J126").Value = 59000
$ws.Range("L126").Value = 59000
$ws.Range("N126").Value = -68880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 500
$ws.Range("I11").Value = 500
$ws.Range("K11").Value = 500
$ws.Range("M11").Value = -360

$ws.Range("H16").Value = 784.75
$ws.Range("I16").Value = 756.8570999999999
$ws.Range("K16").Value = 756.8570999999999
$ws.Range("M16").Value = -469.8570999999999

$ws.Range("H31").Value = 2512.2395
$ws.Range("I31").Value = 931.4878
$ws.Range("J31").Value = 4672.6
$ws.Range("K31").Value = 931.4878
$ws.Range("L31").Value = 4672.6
$ws.Range("M31").Value = -636.4878
$ws.Range("N31").Value = -5262.6

$ws.Range("H33").Value = 1863.1428
$ws.Range("I33").Value = 1863.1428
$ws.Range("K33").Value = 1863.1428
$ws.Range("M33").Value = -1484.1428

$ws.Range("H34").Value = 2512.2395
$ws.Range("I34").Value = 931.4878
$ws.Range("J34").Value = 4672.6
$ws.Range("K34").Value = 931.4878
$ws.Range("L34").Value = 4672.6
$ws.Range("M34").Value = -729.4878
$ws.Range("N34").Value = -5076.6

$ws.Range("H58").Value = 8582.875
$ws.Range("I58").Value = 12556.889
$ws.Range("J58").Value = 3473.4285
$ws.Range("K58").Value = 12556.889
$ws.Range("L58").Value = 3473.4285
$ws.Range("M58").Value = -12353.889
$ws.Range("N58").Value = -3879.4285

$ws.Range("H94").Value = 4176.4165
$ws.Range("I94").Value = 792.2222
$ws.Range("K94").Value = 792.2222
$ws.Range("M94").Value = -341.2222

$ws.Range("H113").Value = 784.75
$ws.Range("I113").Value = 756.8570999999999
$ws.Range("K113").Value = 756.8570999999999
$ws.Range("M113").Value = 1413.1429

$ws.Range("H132").Value = 13891136
$ws.Range("I132").Value = 22728796
$ws.Range("J132").Value = 3384.1428
$ws.Range("K132").Value = 68186388
$ws.Range("L132").Value = 10152.4284
$ws.Range("M132").Value = -68183858
$ws.Range("N132").Value = -15212.4284

$ws.Range("H136").Value = 8582.875
$ws.Range("I136").Value = 12556.889
$ws.Range("J136").Value = 3473.4285
$ws.Range("K136").Value = 37670.667
$ws.Range("L136").Value = 10420.2855
$ws.Range("M136").Value = -35120.667
$ws.Range("N136").Value = -15520.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1467.2439
$ws.Range("I68").Value = 1460
$ws.Range("J68").Value = 1468.25
$ws.Range("K68").Value = 4380
$ws.Range("L68").Value = 4404.75
$ws.Range("M68").Value = -3569
$ws.Range("N68").Value = -6026.75

$ws.Range("H71").Value = 1467.2439
$ws.Range("I71").Value = 1460
$ws.Range("J71").Value = 1468.25
$ws.Range("K71").Value = 13140
$ws.Range("L71").Value = 13214.25
$ws.Range("M71").Value = -9084
$ws.Range("N71").Value = -21326.25

$ws.Range("H107").Value = 1117.84
$ws.Range("I107").Value = 296.15384
$ws.Range("J107").Value = 2008
$ws.Range("K107").Value = 888.4615200000001
$ws.Range("L107").Value = 6024
$ws.Range("M107").Value = 1031.53848
$ws.Range("N107").Value = -9864

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 528
$ws.Range("I6").Value = 528
$ws.Range("K6").Value = 528
$ws.Range("M6").Value = -415

$ws.Range("H16").Value = 528
$ws.Range("I16").Value = 528
$ws.Range("K16").Value = 528
$ws.Range("M16").Value = -278

$ws.Range("H25").Value = 69009
$ws.Range("J25").Value = 69009
$ws.Range("L25").Value = 69009
$ws.Range("N25").Value = -70067

$ws.Range("H102").Value = 3342.8215
$ws.Range("I102").Value = 1409.5555
$ws.Range("J102").Value = 6822.7
$ws.Range("K102").Value = 1409.5555
$ws.Range("L102").Value = 6822.7
$ws.Range("M102").Value = 212.4445000000001
$ws.Range("N102").Value = -10066.7

$ws.Range("H126").Value = 2148.2856
$ws.Range("I126").Value = 1812.4706
$ws.Range("J126").Value = 3575.5
$ws.Range("K126").Value = 5437.4118
$ws.Range("L126").Value = 10726.5
$ws.Range("M126").Value = -2967.4118
$ws.Range("N126").Value = -15666.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 3580
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H15").Value = 3580
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H19").Value = 1900
$ws.Range("J19").Value = 1900
$ws.Range("L19").Value = 1900
$ws.Range("N19").Value = -2240

$ws.Range("H132").Value = 6295972
$ws.Range("I132").Value = 2457.7097
$ws.Range("J132").Value = 15164106
$ws.Range("K132").Value = 7373.1291
$ws.Range("L132").Value = 45492318
$ws.Range("M132").Value = -4843.1291
$ws.Range("N132").Value = -45497378

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 6650.2
$ws.Range("J61").Value = 11350
$ws.Range("L61").Value = 11350
$ws.Range("N61").Value = -11934

$ws.Range("H81").Value = 3742.2222
$ws.Range("I81").Value = 990
$ws.Range("J81").Value = 4528.5713
$ws.Range("K81").Value = 1980
$ws.Range("L81").Value = 9057.142599999999
$ws.Range("M81").Value = -919
$ws.Range("N81").Value = -11179.1426

$ws.Range("H84").Value = 3742.2222
$ws.Range("I84").Value = 990
$ws.Range("J84").Value = 4528.5713
$ws.Range("K84").Value = 9900
$ws.Range("L84").Value = 45285.713
$ws.Range("M84").Value = -4596
$ws.Range("N84").Value = -55893.713

$ws.Range("H98").Value = 49000
$ws.Range("J98").Value = 49000
$ws.Range("L98").Value = 49000
$ws.Range("N98").Value = -54990

$ws.Range("H107").Value = 806.8
$ws.Range("I107").Value = 765
$ws.Range("J107").Value = 834.6667
$ws.Range("K107").Value = 2295
$ws.Range("L107").Value = 2504.0001
$ws.Range("M107").Value = -375
$ws.Range("N107").Value = -6344.0001

$ws.Range("H113").Value = 375.41666
$ws.Range("I113").Value = 356
$ws.Range("J113").Value = 402.6
$ws.Range("K113").Value = 1068
$ws.Range("L113").Value = 1207.8
$ws.Range("M113").Value = 1102
$ws.Range("N113").Value = -5547.8
